$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 2021.2
$ws.Range("I129").Value = 1051
$ws.Range("J129").Value = 2437
$ws.Range("K129").Value = 3153
$ws.Range("L129").Value = 7311
$ws.Range("M129").Value = 1847
$ws.Range("N129").Value = -17311

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3999.9
$ws.Range("I61").Value = 3999.9
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 3999.9
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -3787.9

$ws.Range("H97").Value = 994.6316
$ws.Range("I97").Value = 745.3077
$ws.Range("J97").Value = 1534.8334
$ws.Range("K97").Value = 745.3077
$ws.Range("L97").Value = 1534.8334
$ws.Range("M97").Value = -249.3077
$ws.Range("N97").Value = -2526.8334

$ws.Range("H122").Value = 2218.1428
$ws.Range("I122").Value = 1003
$ws.Range("J122").Value = 3838.3333
$ws.Range("K122").Value = 3009
$ws.Range("L122").Value = 11514.9999
$ws.Range("M122").Value = -559
$ws.Range("N122").Value = -16414.9999

$ws.Range("H136").Value = 3999.9
$ws.Range("I136").Value = 3999.9
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 11999.7
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -9449.700000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2494.1428
$ws.Range("I20").Value = 2494.1428
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 2494.1428
$ws.Range("L20").Value = 0
$ws.Range("M20").Value = -2247.1428

$ws.Range("H26").Value = 42497.5
$ws.Range("I26").Value = 43490.332
$ws.Range("J26").Value = 39519
$ws.Range("K26").Value = 43490.332
$ws.Range("L26").Value = 39519
$ws.Range("M26").Value = -43198.332
$ws.Range("N26").Value = -40103

$ws.Range("N57").ClearContents()
$ws.Range("H57").Value = 70709
$ws.Range("I57").Value = 70709
$ws.Range("J57").Value = 0
$ws.Range("K57").Value = 70709
$ws.Range("L57").Value = 0
$ws.Range("M57").Value = -69989

$ws.Range("H60").Value = 94000
$ws.Range("I60").Value = 0
$ws.Range("J60").Value = 94000
$ws.Range("K60").Value = 0
$ws.Range("L60").Value = 94000
$ws.Range("N60").Value = -95198

$ws.Range("H64").Value = 813.1667
$ws.Range("I64").Value = 749
$ws.Range("J64").Value = 877.3333
$ws.Range("K64").Value = 749
$ws.Range("L64").Value = 877.3333
$ws.Range("M64").Value = -524
$ws.Range("N64").Value = -1327.3333

$ws.Range("H67").Value = 813.1667
$ws.Range("I67").Value = 749
$ws.Range("J67").Value = 877.3333
$ws.Range("K67").Value = 749
$ws.Range("L67").Value = 877.3333
$ws.Range("M67").Value = 31
$ws.Range("N67").Value = -2437.3333

$ws.Range("H135").Value = 204999.5
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 204999.5
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 204999.5
$ws.Range("N135").Value = -215139.5

$ws.Range("N136").ClearContents()
$ws.Range("H136").Value = 70709
$ws.Range("I136").Value = 70709
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 70709
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -65609

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1642.2858
$ws.Range("I16").Value = 1666
$ws.Range("J16").Value = 1624.5
$ws.Range("K16").Value = 1666
$ws.Range("L16").Value = 1624.5
$ws.Range("M16").Value = -1379
$ws.Range("N16").Value = -2198.5

$ws.Range("H32").Value = 1061.2142
$ws.Range("I32").Value = 1061.2142
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 1061.2142
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -745.2141999999999

$ws.Range("H58").Value = 4453.778
$ws.Range("I58").Value = 2000
$ws.Range("J58").Value = 7521
$ws.Range("K58").Value = 2000
$ws.Range("L58").Value = 7521
$ws.Range("M58").Value = -1797
$ws.Range("N58").Value = -7927

$ws.Range("H108").Value = 61296
$ws.Range("I108").Value = 0
$ws.Range("J108").Value = 61296
$ws.Range("K108").Value = 0
$ws.Range("L108").Value = 61296
$ws.Range("N108").Value = -68976

$ws.Range("H113").Value = 1642.2858
$ws.Range("I113").Value = 1666
$ws.Range("J113").Value = 1624.5
$ws.Range("K113").Value = 1666
$ws.Range("L113").Value = 1624.5
$ws.Range("M113").Value = 504
$ws.Range("N113").Value = -5964.5

$ws.Range("H136").Value = 4453.778
$ws.Range("I136").Value = 2000
$ws.Range("J136").Value = 7521
$ws.Range("K136").Value = 6000
$ws.Range("L136").Value = 22563
$ws.Range("M136").Value = -3450
$ws.Range("N136").Value = -27663

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H57").Value = 2490
$ws.Range("I57").Value = 0
$ws.Range("J57").Value = 2490
$ws.Range("K57").Value = 0
$ws.Range("L57").Value = 7470
$ws.Range("N57").Value = -8588

$ws.Range("H125").Value = 3749.75
$ws.Range("I125").Value = 4999.5
$ws.Range("J125").Value = 2500
$ws.Range("K125").Value = 14998.5
$ws.Range("L125").Value = 7500
$ws.Range("M125").Value = -10078.5
$ws.Range("N125").Value = -17340

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("M70").ClearContents()
$ws.Range("H70").Value = 10009
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 10009
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 10009
$ws.Range("N70").Value = -10549

$ws.Range("M73").ClearContents()
$ws.Range("H73").Value = 10009
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 10009
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 10009
$ws.Range("N73").Value = -11881

$ws.Range("H113").Value = 9855.111000000001
$ws.Range("I113").Value = 9346.5
$ws.Range("J113").Value = 10000.429
$ws.Range("K113").Value = 9346.5
$ws.Range("L113").Value = 10000.429
$ws.Range("M113").Value = -7176.5
$ws.Range("N113").Value = -14340.429

$ws.Range("H122").Value = 2964.375
$ws.Range("I122").Value = 2245.8572
$ws.Range("J122").Value = 7994
$ws.Range("K122").Value = 6737.571599999999
$ws.Range("L122").Value = 23982
$ws.Range("M122").Value = -4287.571599999999
$ws.Range("N122").Value = -28882

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("N13").ClearContents()
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 0

$ws.Range("H25").Value = 3583.8462
$ws.Range("I25").Value = 2500
$ws.Range("J25").Value = 3909
$ws.Range("K25").Value = 2500
$ws.Range("L25").Value = 3909
$ws.Range("M25").Value = -2270
$ws.Range("N25").Value = -4369

$ws.Range("H61").Value = 4584.6924
$ws.Range("I61").Value = 2514.7144
$ws.Range("J61").Value = 6999.6665
$ws.Range("K61").Value = 2514.7144
$ws.Range("L61").Value = 6999.6665
$ws.Range("M61").Value = -2312.7144
$ws.Range("N61").Value = -7403.6665

$ws.Range("H113").Value = 4584.6924
$ws.Range("I113").Value = 2514.7144
$ws.Range("J113").Value = 6999.6665
$ws.Range("K113").Value = 2514.7144
$ws.Range("L113").Value = 6999.6665
$ws.Range("M113").Value = -344.7143999999998
$ws.Range("N113").Value = -11339.6665

$ws.Range("M122").ClearContents()
$ws.Range("H122").Value = 3500
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 3500
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 10500
$ws.Range("N122").Value = -15400

$ws.Range("H123").Value = 82000
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 82000
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 82000
$ws.Range("N123").Value = -91800

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 15166667
$ws.Range("I5").Value = 17500000
$ws.Range("J5").Value = 10500000
$ws.Range("K5").Value = 17500000
$ws.Range("L5").Value = 10500000
$ws.Range("M5").Value = -17499888
$ws.Range("N5").Value = -10500224
